# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4" (same column layout/styles),
#        inserted right before "总计" so the tab order matches. ---
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($wb.Worksheets.Item("总计"))

# Re-fetch by name: sheet positions shifted after the Copy(), so any handle obtained
# before the structural change is stale.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template ("2021-Q4") only has 21 data rows; 2022-Q1 needs 22, so insert one more
# data row (with the same row style as the rest) before filling in real values.
$newSheet.Rows.Item(23).Insert()
$newSheet.Range("A22").Copy()
$newSheet.Range("A23").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

# --- 2. Fill the 2022-Q1 fund-holding data ---
$fundData = @(
    @("000362", "国泰聚信价值优势灵活配置混合A", "56.15", "89.00", "6.12", "3.4364", 1),
    @("008415", "国泰大制造两年持有期混合", "23.19", "92.05", "6.29", "1.4587", 1),
    @("012173", "国泰兴泽优选一年持有期混合A", "16.89", "89.83", "6.51", "1.0995", 1),
    @("000363", "国泰聚信价值优势灵活配置混合C", "17.09", "89.00", "6.12", "1.0459", 1),
    @("020010", "国泰金牛创新混合", "16.99", "84.21", "6.06", "1.0296", 2),
    @("001579", "国泰大农业股票", "11.72", "90.32", "5.96", "0.6985", 1),
    @("288001", "华夏经典配置混合", "18.49", "63.85", "3.75", "0.6934", 9),
    @("007835", "国泰鑫睿混合", "9.37", "78.94", "6.36", "0.5959", 1),
    @("012174", "国泰兴泽优选一年持有期混合C", "7.14", "89.83", "6.51", "0.4648", 1),
    @("005244", "国泰聚优价值灵活配置混合A", "7.72", "83.97", "5.36", "0.4138", 2),
    @("161729", "招商 3 年封闭运作瑞利灵活配置混合型", "5.33", "86.78", "4.67", "0.2489", 4),
    @("260112", "景顺长城能源基建混合", "16.49", "60.89", "1.48", "0.2441", 9),
    @("005245", "国泰聚优价值灵活配置混合C", "4.52", "83.97", "5.36", "0.2423", 2),
    @("005945", "工银瑞信可转债优选债券A", "6.62", "46.72", "2.14", "0.1417", 5),
    @("001110", "中欧瑾泉灵活配置混合 - A", "7.39", "22.25", "1.35", "0.0998", 8),
    @("011743", "华夏兴源稳健一年持有期混合型证券投资基金A", "17.89", "23.29", "0.55", "0.0984", 5),
    @("014125", "华夏中证1000指数增强A", "7.03", "89.75", "0.83", "0.0583", 5),
    @("014126", "华夏中证1000指数增强C", "6.09", "89.75", "0.83", "0.0505", 5),
    @("001111", "中欧瑾泉灵活配置混合 - C", "2.20", "22.25", "1.35", "0.0297", 8),
    @("011744", "华夏兴源稳健一年持有期混合型证券投资基金C", "4.45", "23.29", "0.55", "0.0245", 5),
    @("005946", "工银瑞信可转债优选债券C", "0.86", "46.72", "2.14", "0.0184", 5),
    @("005128", "华夏永康添福混合", "1.47", "24.52", "0.55", "0.0081", 9)
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 2).ClearFormats()
    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 3).ClearFormats()
    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 4).ClearFormats()
    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 5).ClearFormats()
    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 6).ClearFormats()
    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 7).ClearFormats()
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# --- 3. Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q1 summary,
#        pushing the existing quarters down (matching the existing un-styled B/C/D pattern). ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 22
$totalSheet.Range("D2").Value = 12.2
